# "made Point use smart-pointer-to-impl"
# Adds a "Points On Line" benchmark column to the existing "No opimizations:"
# results table, and adds a second results table ("Smart Pointer to
# Implementation") with the same set of benchmark columns to its right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Left table ("No opimizations:") - add the "Points On Line" column (F)
# and fill in the previously-missing data points.
# ---------------------------------------------------------------------

$ws.Range("F2").Value = "Points On Line"

$ws.Range("D3").Value = 173
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 1004

$ws.Range("D4").Value = 163
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 1208

$ws.Range("D5").Value = 163
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 981

$ws.Range("D6").Value = 155
$ws.Range("E6").Value = 23
$ws.Range("F6").Value = 950

$ws.Range("C7").Value = 1220
$ws.Range("D7").Value = 172
$ws.Range("E7").Value = 24
$ws.Range("F7").Value = 940

# Extend the merged header cell from B1:E1 to B1:F1
$ws.Range("B1:E1").UnMerge()
$ws.Range("B1:F1").Merge()

# ---------------------------------------------------------------------
# Right table ("Smart Pointer to Implementation") - same column layout,
# starting at column H (column G left blank as a gap between tables).
# ---------------------------------------------------------------------

$ws.Range("H1").Value = "Smart Pointer to Implementation"
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1:L1").Merge()

$ws.Range("H2").Value = "Two Square Points"
$ws.Range("I2").Value = "Two Points 0.5 Apart"
$ws.Range("J2").Value = "Small Triangle"
$ws.Range("K2").Value = "Moser Spindle"
$ws.Range("L2").Value = "Points On Line"

# ---------------------------------------------------------------------
# Mean row - extend the existing AVERAGE formulas across the new columns
# of both tables (D9:F9 and H9:L9 form one shared formula group, same as
# D9:L9 with the empty G9 gap skipped).
# ---------------------------------------------------------------------

$ws.Range("D9:L9").Formula = "=AVERAGE(D3:D7)"
$ws.Range("G9").ClearContents()

$ws.Range("D9:F9").NumberFormat = "0"
$ws.Range("H9:L9").NumberFormat = "0"

# ---------------------------------------------------------------------
# Column widths - autofit the newly populated / newly added columns.
# ---------------------------------------------------------------------

$ws.Columns("D:F").AutoFit()
$ws.Columns("H:L").AutoFit()

# ---------------------------------------------------------------------
# Restore cursor/selection position as recorded in the saved workbook.
# ---------------------------------------------------------------------

$ws.Range("F16").Select()
